# Convert the fromHTMLBodyString() field (fldChar begin/instrText.../fldChar end)
# into plain literal text runs "{ ... }" while keeping the bookmark in place.
# This mirrors unlinking the field and retyping its code as visible text.

$d = $word.ActiveDocument

# Locate the paragraph that contains the m2doc field (fldChar begin/end + instrText).
$targetPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $candidate = $d.Paragraphs.Item($i)
    if ($candidate.Range.Fields.Count -gt 0) {
        $targetPara = $candidate
    }
}

$newParagraphXml = '<w:p w:rsidR="00C52979" w:rsidRDefault="00C52979" w:rsidP="00F5495F">' +
    '<w:r><w:t>{</w:t></w:r>' +
    '<w:r><w:t>m</w:t></w:r>' +
    '<w:r><w:t>:</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> (''&lt;h2 id="starting-with-</w:t></w:r>' +
    '<w:r><w:t>m2doc"&gt;Starting with '' + self.na</w:t></w:r>' +
    '<w:r><w:t>me + ''&lt;/h2&gt;'').from</w:t></w:r>' +
    '<w:r><w:t>HTML</w:t></w:r>' +
    '<w:bookmarkStart w:id="0" w:name="_GoBack"/>' +
    '<w:bookmarkEnd w:id="0"/>' +
    '<w:r><w:t>Body</w:t></w:r>' +
    '<w:r><w:t>String()</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve">}</w:t></w:r>' +
    '</w:p>'

$packageXml = '<?xml version="1.0" standalone="yes"?>' +
    '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
    '<pkg:xmlData>' +
    '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
    '<w:body>' + $newParagraphXml + '</w:body>' +
    '</w:document>' +
    '</pkg:xmlData></pkg:part></pkg:package>'

$targetPara.Range.InsertXML($packageXml) | Out-Null

Write-Output "field unlinked to literal text"
